$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "11-08-2021"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 140000
$ws.Range("C5").Value = 343000
$ws.Range("D5").Value = 140000
$ws.Range("E5").Value = 82000
$ws.Range("F5").Value = 58000
$ws.Range("G5").Value = 4.8
